# dummy change to check connectivity to repository
#
# Adds a small "Hallo / Erwin" note as a new row below the existing parts
# list, and updates the active selection the same way it was left in the
# authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 26: A26 = "Hallo ", B26 = "Erwin"
$ws.Range("A26").Value = "Hallo "
$ws.Range("B26").Value = "Erwin"

# Restore the selection state left behind in the source workbook.
$ws.Range("C30").Select() | Out-Null
